$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.330.28'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.683.81'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.43%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '221.04'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.522'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '30.02'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.05%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0624'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.83%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.926.92'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.58%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.72'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +16.86%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +8.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.679.25'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.00'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.17%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.374.10'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '65.79'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '246.58'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0720'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.997'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.29'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.18'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +5.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.19'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.74'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.22%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.14%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.31%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.27%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.31'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.506.50'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +5.50%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.23%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +5.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.75'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.587'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '78.90'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +10.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.32'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.851'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0504'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.84%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.997'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.31%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '51.93'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -6.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.817.76'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.78%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.42'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '95.35'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +6.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0113'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.85%  '
